# demo/simple_demo.xlsx - "updated demo" commit
# - Updated the selenium demo HTML path used by __main__ imports (q1 step),
#   referenced from both the tc0001 and tc0002 test-case sheets.
# - The column holding that longer path (column G) was widened on the two
#   sheets where it now shows the longer text, and the active-cell selection
#   on those two sheets moved to reflect where the user was working.

$wb = $excel.ActiveWorkbook

$wsTr0001 = $wb.Worksheets.Item("tr0001")
$wsTc0001 = $wb.Worksheets.Item("tc0001")
$wsTc0002 = $wb.Worksheets.Item("tc0002")

$newPath = "demo/demo-html/page1.html"

# --- tc0001 (sheet2): update the test-page path and widen column G ---
$wsTc0001.Range("G3").Value = $newPath
$wsTc0001.Columns.Item(7).ColumnWidth = 28.8333333333333

# --- tc0002 (sheet3): same path update, used twice on this sheet ---
$wsTc0002.Range("G3").Value = $newPath
$wsTc0002.Range("G8").Value = $newPath
$wsTc0002.Columns.Item(7).ColumnWidth = 38

# --- update the saved selection / active cell on each changed sheet ---
[void]$wsTc0001.Activate()
$wsTc0001.Range("G3").Select() | Out-Null

[void]$wsTc0002.Activate()
$wsTc0002.Range("G16").Select() | Out-Null
